$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row below mirrors one changed row from the source diff.
# Column D (Price) values are forced to Text storage (NumberFormat "@")
# so numeric-looking strings like "318.22" are not reinterpreted as
# numbers by Excel, then the style is reset to "Normal" so no stray
# cell formatting is left behind.

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value2 = "43.074.81"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Range("E2").Value = "  +0.60%  "

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value2 = "2.545.37"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Range("E3").Value = "  +0.55%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value2 = "318.22"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Range("E5").Value = "  +4.40%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value2 = "96.26"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Range("E6").Value = "  -2.46%  "

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value2 = "0.581"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Range("E7").Value = "  -0.34%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value2 = "0.534"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Range("E9").Value = "  -2.08%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value2 = "36.55"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Range("E10").Value = "  -1.02%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value2 = "0.0814"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Range("E11").Value = "  +0.27%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value2 = "7.66"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Range("E12").Value = "  -0.95%  "

$ws.Range("E13").Value = "  +0.66%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value2 = "2.935.83"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Range("E14").Value = "  +0.50%  "

$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value2 = "15.58"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Range("E15").Value = "  +3.22%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value2 = "2.551.13"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Range("E16").Value = "  -0.51%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value2 = "0.855"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Range("E17").Value = "  -1.58%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value2 = "42.896.98"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Range("E18").Value = "  +0.05%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value2 = "13.11"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Range("E19").Value = "  +0.70%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value2 = "6.62"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Range("E20").Value = "  +2.11%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value2 = "0.0₃0970"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Range("E21").Value = "  -1.20%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value2 = "70.55"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Range("E22").Value = "  -1.34%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value2 = "252.78"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Range("E23").Value = "  -0.10%  "

$ws.Range("E24").Value = "  +1.60%  "

$ws.Range("E25").Value = "  -0.79%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value2 = "26.91"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Range("E26").Value = "  -0.30%  "

$ws.Range("E27").Value = "  -0.02%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value2 = "2.42"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Range("E28").Value = "  +4.27%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value2 = "39.89"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Range("E29").Value = "  +3.65%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value2 = "10.23"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Range("E30").Value = "  -2.04%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value2 = "6.10"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Range("E31").Value = "  -0.20%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value2 = "155.19"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Range("E32").Value = "  -1.89%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value2 = "2.15"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Range("E33").Value = "  +2.14%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value2 = "19.18"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Range("E34").Value = "  +5.04%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value2 = "3.34"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Range("E35").Value = "  +0.63%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value2 = "0.0793"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Range("E36").Value = "  +0.06%  "

$ws.Range("E37").Value = "  +0.32%  "

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value2 = "0.113"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Range("E38").Value = "  -2.04%  "

$ws.Range("E39").Value = "  -0.39%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value2 = "23.80"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Range("E40").Value = "  -1.44%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value2 = "2.30"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Range("E41").Value = "  +10.19%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value2 = "3.83"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Range("E42").Value = "  -1.72%  "

$ws.Range("B43").Value = "NEARProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value2 = "3.35"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Range("E43").Value = "  -2.49%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value2 = "0.0304"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Range("E44").Value = "  +0.58%  "

$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value2 = "1.00"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Range("E45").Value = "  +0.17%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value2 = "2.023.75"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Range("E46").Value = "  -0.90%  "

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value2 = "86.15"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Range("E47").Value = "  +0.22%  "

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value2 = "8.87"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Range("E48").Value = "  -1.18%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value2 = "2.788.08"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Range("E49").Value = "  +0.36%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value2 = "74.79"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Range("E50").Value = "  +2.74%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value2 = "103.13"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Range("E51").Value = "  +0.19%  "
